# Update countries & provincias Spain
# This script updates the COVID-19 stats table on the "Pais" sheet:
#  - Updates case numbers for several countries (Estados Unidos, Paises Bajos,
#    Croacia, Niger)
#  - Updates Sri Lanka's stats and re-sorts it into its new position (the
#    table is sorted by column B - "Casos totales" - descending), which
#    pushes Eslovenia and Kenia down one row each
#  - Updates the "Datos actualizados..." timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1746311
$ws.Range("C4").Value = 508
$ws.Range("D4").Value = 490256
$ws.Range("E4").Value = 1153939
$ws.Range("G4").Value = 9
$ws.Range("H4").Value = 102116

# --- Paises Bajos (row 24) ---
$ws.Range("B24").Value = 45950
$ws.Range("C24").Value = 182
$ws.Range("G24").Value = 32
$ws.Range("H24").Value = 5903

# --- Croacia (row 88) ---
$ws.Range("B88").Value = 2245
$ws.Range("C88").Value = 1
$ws.Range("D88").Value = 2051
$ws.Range("E88").Value = 92
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 102

# --- Sri Lanka / Eslovenia / Kenia re-sort (rows 100-102) ---
# Before: row100=Eslovenia, row101=Kenia, row102=Sri Lanka
# After (Sri Lanka updated & moved up, so it now sits right after Nueva
# Zelanda, above Eslovenia and Kenia which both shift down one row):
#   row100=Sri Lanka (new data), row101=Eslovenia (old row100 data),
#   row102=Kenia (old row101 data)
$ws.Range("A101").Value = "Eslovenia"
$ws.Range("B101").Value = 1473
$ws.Range("C101").Value = 2
$ws.Range("D101").Value = 1356
$ws.Range("E101").Value = 9
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 108

$ws.Range("A102").Value = "Kenia"
$ws.Range("B102").Value = 1471
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 408
$ws.Range("E102").Value = 1008
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 55

$ws.Range("A100").Value = "Sri Lanka"
$ws.Range("B100").Value = 1486
$ws.Range("C100").Value = 17
$ws.Range("D100").Value = 745
$ws.Range("E100").Value = 731
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 10

# --- Niger (row 118) ---
$ws.Range("B118").Value = 955
$ws.Range("C118").Value = 3
$ws.Range("E118").Value = 95
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 64

# --- Update "last updated" timestamp (cell A1) ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 14:40"

$wb.Save()
